$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new fly_part2 vial/treatment rows (rows 614-681), matching the
# new batch of vials at time_hours = 291.
$ws.Range("A614").Value2 = 1
$ws.Range("B614").Value2 = "conditioned"
$ws.Range("C614").Value2 = 291
$ws.Range("D614").Value2 = 3
$ws.Range("E614").Value2 = 0
$ws.Range("A615").Value2 = 1
$ws.Range("B615").Value2 = "unconditioned"
$ws.Range("C615").Value2 = 291
$ws.Range("D615").Value2 = 0
$ws.Range("E615").Value2 = 0
$ws.Range("A616").Value2 = 2
$ws.Range("B616").Value2 = "conditioned"
$ws.Range("C616").Value2 = 291
$ws.Range("D616").Value2 = 0
$ws.Range("E616").Value2 = 2
$ws.Range("A617").Value2 = 2
$ws.Range("B617").Value2 = "unconditioned"
$ws.Range("C617").Value2 = 291
$ws.Range("D617").Value2 = 1
$ws.Range("E617").Value2 = 1
$ws.Range("A618").Value2 = 3
$ws.Range("B618").Value2 = "conditioned"
$ws.Range("C618").Value2 = 291
$ws.Range("D618").Value2 = 0
$ws.Range("E618").Value2 = 0
$ws.Range("A619").Value2 = 3
$ws.Range("B619").Value2 = "unconditioned"
$ws.Range("C619").Value2 = 291
$ws.Range("D619").Value2 = 2
$ws.Range("E619").Value2 = 0
$ws.Range("A620").Value2 = 4
$ws.Range("B620").Value2 = "conditioned"
$ws.Range("C620").Value2 = 291
$ws.Range("D620").Value2 = 5
$ws.Range("E620").Value2 = 3
$ws.Range("A621").Value2 = 4
$ws.Range("B621").Value2 = "unconditioned"
$ws.Range("C621").Value2 = 291
$ws.Range("D621").Value2 = 0
$ws.Range("E621").Value2 = 0
$ws.Range("A622").Value2 = 5
$ws.Range("B622").Value2 = "conditioned"
$ws.Range("C622").Value2 = 291
$ws.Range("D622").Value2 = 1
$ws.Range("E622").Value2 = 1
$ws.Range("A623").Value2 = 5
$ws.Range("B623").Value2 = "unconditioned"
$ws.Range("C623").Value2 = 291
$ws.Range("D623").Value2 = 1
$ws.Range("E623").Value2 = 1
$ws.Range("A624").Value2 = 6
$ws.Range("B624").Value2 = "conditioned"
$ws.Range("C624").Value2 = 291
$ws.Range("D624").Value2 = 0
$ws.Range("E624").Value2 = 0
$ws.Range("A625").Value2 = 6
$ws.Range("B625").Value2 = "unconditioned"
$ws.Range("C625").Value2 = 291
$ws.Range("D625").Value2 = 4
$ws.Range("E625").Value2 = 0
$ws.Range("A626").Value2 = 7
$ws.Range("B626").Value2 = "conditioned"
$ws.Range("C626").Value2 = 291
$ws.Range("D626").Value2 = 3
$ws.Range("E626").Value2 = 0
$ws.Range("A627").Value2 = 7
$ws.Range("B627").Value2 = "unconditioned"
$ws.Range("C627").Value2 = 291
$ws.Range("D627").Value2 = 0
$ws.Range("E627").Value2 = 0
$ws.Range("A628").Value2 = 8
$ws.Range("B628").Value2 = "conditioned"
$ws.Range("C628").Value2 = 291
$ws.Range("D628").Value2 = 0
$ws.Range("E628").Value2 = 0
$ws.Range("A629").Value2 = 8
$ws.Range("B629").Value2 = "unconditioned"
$ws.Range("C629").Value2 = 291
$ws.Range("D629").Value2 = 29
$ws.Range("E629").Value2 = 12
$ws.Range("A630").Value2 = 9
$ws.Range("B630").Value2 = "conditioned"
$ws.Range("C630").Value2 = 291
$ws.Range("D630").Value2 = 0
$ws.Range("E630").Value2 = 1
$ws.Range("A631").Value2 = 9
$ws.Range("B631").Value2 = "unconditioned"
$ws.Range("C631").Value2 = 291
$ws.Range("D631").Value2 = 6
$ws.Range("E631").Value2 = 2
$ws.Range("A632").Value2 = 10
$ws.Range("B632").Value2 = "conditioned"
$ws.Range("C632").Value2 = 291
$ws.Range("D632").Value2 = 2
$ws.Range("E632").Value2 = 2
$ws.Range("A633").Value2 = 10
$ws.Range("B633").Value2 = "unconditioned"
$ws.Range("C633").Value2 = 291
$ws.Range("D633").Value2 = 2
$ws.Range("E633").Value2 = 0
$ws.Range("A634").Value2 = 11
$ws.Range("B634").Value2 = "conditioned"
$ws.Range("C634").Value2 = 291
$ws.Range("D634").Value2 = 0
$ws.Range("E634").Value2 = 0
$ws.Range("A635").Value2 = 11
$ws.Range("B635").Value2 = "unconditioned"
$ws.Range("C635").Value2 = 291
$ws.Range("D635").Value2 = 0
$ws.Range("E635").Value2 = 0
$ws.Range("A636").Value2 = 12
$ws.Range("B636").Value2 = "conditioned"
$ws.Range("C636").Value2 = 291
$ws.Range("D636").Value2 = 0
$ws.Range("E636").Value2 = 0
$ws.Range("A637").Value2 = 12
$ws.Range("B637").Value2 = "unconditioned"
$ws.Range("C637").Value2 = 291
$ws.Range("D637").Value2 = 0
$ws.Range("E637").Value2 = 1
$ws.Range("A638").Value2 = 13
$ws.Range("B638").Value2 = "conditioned"
$ws.Range("C638").Value2 = 291
$ws.Range("D638").Value2 = 1
$ws.Range("E638").Value2 = 4
$ws.Range("A639").Value2 = 13
$ws.Range("B639").Value2 = "unconditioned"
$ws.Range("C639").Value2 = 291
$ws.Range("D639").Value2 = 3
$ws.Range("E639").Value2 = 4
$ws.Range("A640").Value2 = 14
$ws.Range("B640").Value2 = "conditioned"
$ws.Range("C640").Value2 = 291
$ws.Range("D640").Value2 = 1
$ws.Range("E640").Value2 = 1
$ws.Range("A641").Value2 = 14
$ws.Range("B641").Value2 = "unconditioned"
$ws.Range("C641").Value2 = 291
$ws.Range("D641").Value2 = 1
$ws.Range("E641").Value2 = 5
$ws.Range("A642").Value2 = 15
$ws.Range("B642").Value2 = "conditioned"
$ws.Range("C642").Value2 = 291
$ws.Range("D642").Value2 = 1
$ws.Range("E642").Value2 = 2
$ws.Range("A643").Value2 = 15
$ws.Range("B643").Value2 = "unconditioned"
$ws.Range("C643").Value2 = 291
$ws.Range("D643").Value2 = 0
$ws.Range("E643").Value2 = 0
$ws.Range("A644").Value2 = 16
$ws.Range("B644").Value2 = "conditioned"
$ws.Range("C644").Value2 = 291
$ws.Range("D644").Value2 = 7
$ws.Range("E644").Value2 = 5
$ws.Range("A645").Value2 = 16
$ws.Range("B645").Value2 = "unconditioned"
$ws.Range("C645").Value2 = 291
$ws.Range("D645").Value2 = 0
$ws.Range("E645").Value2 = 0
$ws.Range("A646").Value2 = 17
$ws.Range("B646").Value2 = "conditioned"
$ws.Range("C646").Value2 = 291
$ws.Range("D646").Value2 = 0
$ws.Range("E646").Value2 = 1
$ws.Range("A647").Value2 = 17
$ws.Range("B647").Value2 = "unconditioned"
$ws.Range("C647").Value2 = 291
$ws.Range("D647").Value2 = 0
$ws.Range("E647").Value2 = 0
$ws.Range("A648").Value2 = 18
$ws.Range("B648").Value2 = "conditioned"
$ws.Range("C648").Value2 = 291
$ws.Range("D648").Value2 = 1
$ws.Range("E648").Value2 = 1
$ws.Range("A649").Value2 = 18
$ws.Range("B649").Value2 = "unconditioned"
$ws.Range("C649").Value2 = 291
$ws.Range("D649").Value2 = 3
$ws.Range("E649").Value2 = 2
$ws.Range("A650").Value2 = 19
$ws.Range("B650").Value2 = "conditioned"
$ws.Range("C650").Value2 = 291
$ws.Range("D650").Value2 = 0
$ws.Range("E650").Value2 = 0
$ws.Range("A651").Value2 = 19
$ws.Range("B651").Value2 = "unconditioned"
$ws.Range("C651").Value2 = 291
$ws.Range("D651").Value2 = 4
$ws.Range("E651").Value2 = 2
$ws.Range("A652").Value2 = 20
$ws.Range("B652").Value2 = "conditioned"
$ws.Range("C652").Value2 = 291
$ws.Range("D652").Value2 = 0
$ws.Range("E652").Value2 = 0
$ws.Range("A653").Value2 = 20
$ws.Range("B653").Value2 = "unconditioned"
$ws.Range("C653").Value2 = 291
$ws.Range("D653").Value2 = 2
$ws.Range("E653").Value2 = 1
$ws.Range("A654").Value2 = 21
$ws.Range("B654").Value2 = "conditioned"
$ws.Range("C654").Value2 = 291
$ws.Range("D654").Value2 = 0
$ws.Range("E654").Value2 = 1
$ws.Range("A655").Value2 = 21
$ws.Range("B655").Value2 = "unconditioned"
$ws.Range("C655").Value2 = 291
$ws.Range("D655").Value2 = 2
$ws.Range("E655").Value2 = 1
$ws.Range("A656").Value2 = 22
$ws.Range("B656").Value2 = "conditioned"
$ws.Range("C656").Value2 = 291
$ws.Range("D656").Value2 = 0
$ws.Range("E656").Value2 = 0
$ws.Range("A657").Value2 = 22
$ws.Range("B657").Value2 = "unconditioned"
$ws.Range("C657").Value2 = 291
$ws.Range("D657").Value2 = 4
$ws.Range("E657").Value2 = 0
$ws.Range("A658").Value2 = 23
$ws.Range("B658").Value2 = "conditioned"
$ws.Range("C658").Value2 = 291
$ws.Range("D658").Value2 = 0
$ws.Range("E658").Value2 = 0
$ws.Range("A659").Value2 = 23
$ws.Range("B659").Value2 = "unconditioned"
$ws.Range("C659").Value2 = 291
$ws.Range("D659").Value2 = 0
$ws.Range("E659").Value2 = 0
$ws.Range("A660").Value2 = 24
$ws.Range("B660").Value2 = "conditioned"
$ws.Range("C660").Value2 = 291
$ws.Range("D660").Value2 = 0
$ws.Range("E660").Value2 = 0
$ws.Range("A661").Value2 = 24
$ws.Range("B661").Value2 = "unconditioned"
$ws.Range("C661").Value2 = 291
$ws.Range("D661").Value2 = 2
$ws.Range("E661").Value2 = 1
$ws.Range("A662").Value2 = 25
$ws.Range("B662").Value2 = "conditioned"
$ws.Range("C662").Value2 = 291
$ws.Range("D662").Value2 = 3
$ws.Range("E662").Value2 = 0
$ws.Range("A663").Value2 = 25
$ws.Range("B663").Value2 = "unconditioned"
$ws.Range("C663").Value2 = 291
$ws.Range("D663").Value2 = 5
$ws.Range("E663").Value2 = 4
$ws.Range("A664").Value2 = 26
$ws.Range("B664").Value2 = "conditioned"
$ws.Range("C664").Value2 = 291
$ws.Range("D664").Value2 = 9
$ws.Range("E664").Value2 = 6
$ws.Range("A665").Value2 = 26
$ws.Range("B665").Value2 = "unconditioned"
$ws.Range("C665").Value2 = 291
$ws.Range("D665").Value2 = 4
$ws.Range("E665").Value2 = 0
$ws.Range("A666").Value2 = 27
$ws.Range("B666").Value2 = "conditioned"
$ws.Range("C666").Value2 = 291
$ws.Range("D666").Value2 = 1
$ws.Range("E666").Value2 = 2
$ws.Range("A667").Value2 = 27
$ws.Range("B667").Value2 = "unconditioned"
$ws.Range("C667").Value2 = 291
$ws.Range("D667").Value2 = 8
$ws.Range("E667").Value2 = 2
$ws.Range("A668").Value2 = 28
$ws.Range("B668").Value2 = "conditioned"
$ws.Range("C668").Value2 = 291
$ws.Range("D668").Value2 = 0
$ws.Range("E668").Value2 = 6
$ws.Range("A669").Value2 = 28
$ws.Range("B669").Value2 = "unconditioned"
$ws.Range("C669").Value2 = 291
$ws.Range("D669").Value2 = 5
$ws.Range("E669").Value2 = 1
$ws.Range("A670").Value2 = 29
$ws.Range("B670").Value2 = "conditioned"
$ws.Range("C670").Value2 = 291
$ws.Range("D670").Value2 = 1
$ws.Range("E670").Value2 = 1
$ws.Range("A671").Value2 = 29
$ws.Range("B671").Value2 = "unconditioned"
$ws.Range("C671").Value2 = 291
$ws.Range("D671").Value2 = 1
$ws.Range("E671").Value2 = 1
$ws.Range("A672").Value2 = 30
$ws.Range("B672").Value2 = "conditioned"
$ws.Range("C672").Value2 = 291
$ws.Range("D672").Value2 = 1
$ws.Range("E672").Value2 = 2
$ws.Range("A673").Value2 = 30
$ws.Range("B673").Value2 = "unconditioned"
$ws.Range("C673").Value2 = 291
$ws.Range("D673").Value2 = 2
$ws.Range("E673").Value2 = 4
$ws.Range("A674").Value2 = 31
$ws.Range("B674").Value2 = "conditioned"
$ws.Range("C674").Value2 = 291
$ws.Range("D674").Value2 = 1
$ws.Range("E674").Value2 = 1
$ws.Range("A675").Value2 = 31
$ws.Range("B675").Value2 = "unconditioned"
$ws.Range("C675").Value2 = 291
$ws.Range("D675").Value2 = 1
$ws.Range("E675").Value2 = 1
$ws.Range("A676").Value2 = 32
$ws.Range("B676").Value2 = "conditioned"
$ws.Range("C676").Value2 = 291
$ws.Range("D676").Value2 = 3
$ws.Range("E676").Value2 = 2
$ws.Range("A677").Value2 = 32
$ws.Range("B677").Value2 = "unconditioned"
$ws.Range("C677").Value2 = 291
$ws.Range("D677").Value2 = 3
$ws.Range("E677").Value2 = 0
$ws.Range("A678").Value2 = 33
$ws.Range("B678").Value2 = "conditioned"
$ws.Range("C678").Value2 = 291
$ws.Range("D678").Value2 = 1
$ws.Range("E678").Value2 = 1
$ws.Range("A679").Value2 = 33
$ws.Range("B679").Value2 = "unconditioned"
$ws.Range("C679").Value2 = 291
$ws.Range("D679").Value2 = 1
$ws.Range("E679").Value2 = 1
$ws.Range("A680").Value2 = 34
$ws.Range("B680").Value2 = "conditioned"
$ws.Range("C680").Value2 = 291
$ws.Range("D680").Value2 = 1
$ws.Range("E680").Value2 = 1
$ws.Range("A681").Value2 = 34
$ws.Range("B681").Value2 = "unconditioned"
$ws.Range("C681").Value2 = 291
$ws.Range("D681").Value2 = 1
$ws.Range("E681").Value2 = 1

# Reflect the author's updated view state: new selection near the bottom of
# the appended data and the zoom level used while reviewing it.
$win = $excel.ActiveWindow
$win.Zoom = 125
$ws.Range("F675").Select() | Out-Null
